$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Update data value (O2): "8417" -> "8419"
$ws.Range("O2").Value = "8419"

# Update header label (N1): "descripcion" -> "descripcionFactura"
$ws.Range("N1").Value = "descripcionFactura"

# Update the active selection on the sheet to O1
$ws.Range("O1").Select()
